$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the invalid RAM (memory) constant values for the first two VM sizes
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2

# Re-apply the default font across the whole sheet (columns A:XFD), which
# Excel records as an explicit style even though it matches the existing look
$ws.Cells.Font.Name = "Calibri"
$ws.Cells.Font.Size = 11

# Update the active selection left by the editing session
$ws.Range("D5").Select() | Out-Null
